# Update the date heading.
$d = $word.ActiveDocument
$d.Content.Find.Execute("2025-11-03 Monday", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "2025-11-04 Tuesday", 2)

# Update the division-fact answers in the single results table.
# The table has 20 rows x 5 columns; data lives in rows 1, 5, 9, 13, 17
# (the other rows are spacer/blank rows).
$t = $d.Tables.Item(1)

$updates = @(
    @{Row=1;  Col=1; New="713÷3=237, 2"},
    @{Row=1;  Col=2; New="574÷3=191, 1"},
    @{Row=1;  Col=3; New="136÷2=68, 0"},
    @{Row=1;  Col=4; New="374÷4=93, 2"},
    @{Row=1;  Col=5; New="455÷5=91, 0"},

    @{Row=5;  Col=1; New="699÷9=77, 6"},
    @{Row=5;  Col=2; New="179÷7=25, 4"},
    @{Row=5;  Col=3; New="154÷6=25, 4"},
    @{Row=5;  Col=4; New="867÷7=123, 6"},
    @{Row=5;  Col=5; New="853÷6=142, 1"},

    @{Row=9;  Col=1; New="532÷6=88, 4"},
    @{Row=9;  Col=2; New="432÷8=54, 0"},
    @{Row=9;  Col=3; New="958÷6=159, 4"},
    @{Row=9;  Col=4; New="736÷6=122, 4"},
    @{Row=9;  Col=5; New="995÷6=165, 5"},

    @{Row=13; Col=1; New="606÷2=303, 0"},
    @{Row=13; Col=2; New="722÷4=180, 2"},
    @{Row=13; Col=3; New="335÷6=55, 5"},
    @{Row=13; Col=4; New="135÷9=15, 0"},
    @{Row=13; Col=5; New="717÷6=119, 3"},

    @{Row=17; Col=1; New="173÷9=19, 2"},
    @{Row=17; Col=2; New="898÷6=149, 4"},
    @{Row=17; Col=3; New="935÷7=133, 4"},
    @{Row=17; Col=4; New="619÷6=103, 1"},
    @{Row=17; Col=5; New="347÷5=69, 2"}
)

foreach ($u in $updates) {
    $t.Cell($u.Row, $u.Col).Range.Text = $u.New
}
